$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume(1h) figures per the latest symbol list refresh.

$ws.Range("D2").Value = "'300.55"
$ws.Range("E2").Value = "'-4.77%"
$ws.Range("D3").Value = "'35.18"
$ws.Range("E3").Value = "'-0.50%"
$ws.Range("E4").Value = "'-0.88%"
$ws.Range("D5").Value = "'0.07941"
$ws.Range("E5").Value = "'-2.97%"
$ws.Range("D6").Value = "'1.900"
$ws.Range("E6").Value = "'-7.90%"
$ws.Range("D7").Value = "'7.778"
$ws.Range("E7").Value = "'-2.07%"
$ws.Range("E8").Value = "'-2.64%"
$ws.Range("D9").Value = "'0.9273"
$ws.Range("E9").Value = "'-0.45%"
$ws.Range("D10").Value = "'0.1334"
$ws.Range("E10").Value = "'29.38%"
$ws.Range("D11").Value = "'0.1896"
$ws.Range("E11").Value = "'-1.34%"
$ws.Range("D12").Value = "'0.09143"
$ws.Range("E12").Value = "'-0.19%"
$ws.Range("E13").Value = "'-3.75%"
$ws.Range("D14").Value = "'0.09902"
$ws.Range("E14").Value = "'-0.11%"
$ws.Range("D15").Value = "'0.001394"
$ws.Range("E15").Value = "'-3.12%"
$ws.Range("D16").Value = "'0.005735"
$ws.Range("E16").Value = "'-0.08%"
$ws.Range("D17").Value = "'3.513"
$ws.Range("E17").Value = "'1.27%"
$ws.Range("D18").Value = "'2.936"
$ws.Range("E18").Value = "'-0.74%"
$ws.Range("E19").Value = "'-0.08%"
$ws.Range("D20").Value = "'0.1293"
$ws.Range("E20").Value = "'-1.24%"
$ws.Range("D21").Value = "'5.038"
$ws.Range("E21").Value = "'-1.22%"
$ws.Range("D22").Value = "'0.2400"
$ws.Range("E22").Value = "'8.46%"
$ws.Range("E23").Value = "'-1.10%"
$ws.Range("D25").Value = "'0.004751"
$ws.Range("E25").Value = "'-0.93%"
$ws.Range("E26").Value = "'-1.62%"
$ws.Range("D27").Value = "'0.0003002"
$ws.Range("E27").Value = "'-33.32%"
$ws.Range("D39").Value = "'0.01887"
$ws.Range("E39").Value = "'-4.93%"
$ws.Range("D40").Value = "'0.04699"
$ws.Range("E40").Value = "'-5.28%"
$ws.Range("D41").Value = "'0.007340"
$ws.Range("E41").Value = "'-3.06%"
$ws.Range("D42").Value = "'0.01004"
$ws.Range("E42").Value = "'27.57%"
$ws.Range("D43").Value = "'0.1318"
$ws.Range("E43").Value = "'-4.81%"
$ws.Range("D44").Value = "'0.002111"
$ws.Range("E44").Value = "'-6.66%"
$ws.Range("D45").Value = "'0.009318"
$ws.Range("E45").Value = "'-20.85%"
$ws.Range("D46").Value = "'0.00006261"
$ws.Range("E46").Value = "'-5.30%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.03%"
$ws.Range("D48").Value = "'64.70"
$ws.Range("E48").Value = "'28.40%"
$ws.Range("D49").Value = "'0.001660"
$ws.Range("E49").Value = "'-2.41%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.03%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.03%"
